$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price/volume snapshot with the latest scrape.
# Prices and percentages are plain text in this sheet (no numeric
# formatting), so a leading apostrophe is used to force Excel to store
# them as literal text instead of auto-converting to a Number/Percent.

# Row 2: BNB
$ws.Range("D2").Value = "'247.03"
$ws.Range("E2").Value = "'0.72%"

# Row 3: OKB
$ws.Range("D3").Value = "'29.31"
$ws.Range("E3").Value = "'7.26%"

# Row 4: HuobiToken
$ws.Range("D4").Value = "'5.188"
$ws.Range("E4").Value = "'1.47%"

# Row 5: Cronos
$ws.Range("D5").Value = "'0.05734"
$ws.Range("E5").Value = "'0.85%"

# Row 6: KuCoinToken
$ws.Range("E6").Value = "'0.37%"

# Row 7: GateToken
$ws.Range("D7").Value = "'3.098"
$ws.Range("E7").Value = "'2.97%"

# Row 8: MXToken
$ws.Range("D8").Value = "'0.8591"
$ws.Range("E8").Value = "'4.80%"

# Row 9: FTXToken
$ws.Range("D9").Value = "'0.8702"
$ws.Range("E9").Value = "'1.08%"

# Row 10: WazirX
$ws.Range("D10").Value = "'0.1366"
$ws.Range("E10").Value = "'2.54%"

# Row 11: MandalaExchangeToken
$ws.Range("D11").Value = "'0.07073"
$ws.Range("E11").Value = "'1.71%"

# Row 12: BitrueCoin
$ws.Range("D12").Value = "'0.03005"
$ws.Range("E12").Value = "'5.34%"

# Row 13: BitMartToken
$ws.Range("D13").Value = "'0.09384"
$ws.Range("E13").Value = "'-0.04%"

# Row 14: BitForexToken
$ws.Range("D14").Value = "'0.001525"
$ws.Range("E14").Value = "'-0.09%"

# Row 15: CoinExToken
$ws.Range("D15").Value = "'0.04141"
$ws.Range("E15").Value = "'1.73%"

# Row 16: One
$ws.Range("D16").Value = "'0.0005988"
$ws.Range("E16").Value = "'0.16%"

# Row 17: TigerCash
$ws.Range("D17").Value = "'0.006023"
$ws.Range("E17").Value = "'-3.08%"

# Row 18: UpBots
$ws.Range("E18").Value = "'5,224.98%"

# Row 19: LEO
$ws.Range("D19").Value = "'3.491"

# Row 20: BTSEToken
$ws.Range("D20").Value = "'2.283"
$ws.Range("E20").Value = "'-1.47%"

# Row 21: BitpandaEcosystemToken
$ws.Range("E21").Value = "'0.64%"

# Row 22: LiechtensteinCryptoassetsExchange
$ws.Range("D22").Value = "'0.03381"
$ws.Range("E22").Value = "'4.95%"

# Row 23: ProBitToken
$ws.Range("D23").Value = "'0.1286"
$ws.Range("E23").Value = "'0.94%"

# Row 24: MCDex
$ws.Range("D24").Value = "'3.464"
$ws.Range("E24").Value = "'-2.49%"

# Row 26: HotbitToken
$ws.Range("E26").Value = "'12.05%"

# Row 27: BitKan
$ws.Range("D27").Value = "'0.001225"
$ws.Range("E27").Value = "'0.69%"

# Row 28: NitroEx
$ws.Range("D28").Value = "'0.0001209"
$ws.Range("E28").Value = "'2.49%"

# Row 40: IDEX
$ws.Range("D40").Value = "'0.03751"
$ws.Range("E40").Value = "'0.75%"

# Row 41: BKEXToken
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.005747"
$ws.Range("E41").Value = "'-2.92%"

# Row 42: CEJI
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1072"
$ws.Range("E42").Value = "'1.35%"

# Row 43: KickToken
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002427"
$ws.Range("E43").Value = "'34.86%"

# Row 44: LocalTraders
$ws.Range("D44").Value = "'0.009417"
$ws.Range("E44").Value = "'-3.06%"

# Row 45: CoinLion
$ws.Range("D45").Value = "'0.00005261"
$ws.Range("E45").Value = "'3.13%"

# Row 46: Kangarootoken
$ws.Range("E46").Value = "'-0.04%"

# Row 47: CoinbaseStockToken
$ws.Range("E47").Value = "'-43.58%"

# Row 48: BOLO
$ws.Range("E48").Value = "'-10.01%"

# Row 49: CryptobidCoin
$ws.Range("E49").Value = "'-0.04%"

# Row 50: SpecialPowerGold
$ws.Range("E50").Value = "'-0.04%"
